$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns("D:D").Insert()

# Copy number formatting from column E (the old column D, now shifted) into new column D
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate new column D with the latest reporting period data
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 6319100
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 1700
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 4202200
$ws.Range("D18").Value = 2116900
$ws.Range("D20").Value = 147000
$ws.Range("D21").Value = 2340300
$ws.Range("D22").Value = 48700
$ws.Range("D23").Value = 2215200
$ws.Range("D24").Value = 503700
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 1711500
$ws.Range("D27").Value = 1715600
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = -968800
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -147000
$ws.Range("D33").Value = 746800
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 746800
$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 6910600
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 847900
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 3535900
$ws.Range("D48").Value = 535000
$ws.Range("D49").Value = 2333400
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = "NA"
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 14383500
$ws.Range("D57").Value = 226900
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 1568100
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 728500
$ws.Range("D62").Value = 126500
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 4484300
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 10217900
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 9899200
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = 746800
$ws.Range("D83").Value = 76400
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 2229700
$ws.Range("D91").Value = -106500
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -290400
$ws.Range("D96").Value = -2116900
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -3761700
$ws.Range("D101").Value = -16700
$ws.Range("D102").Value = -1839100

# Correct rows whose historical figures were also refreshed (not a pure column shift)
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "NA"
$ws.Range("H29").Value = "NA"
$ws.Range("I29").Value = "NA"
$ws.Range("J29").Value = "NA"
$ws.Range("K29").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 34600
$ws.Range("E58").Value = 0
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 84500
